# Add new worksheet "Acceptance test table clear" at the end of the workbook
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Acceptance test table clear"

# Header row (reuses existing shared strings: Test ID / Description / Expected Result / Precondition / Comments)
$ws.Range("A1").Value = "Test ID"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Expected Result"
$ws.Range("D1").Value = "Precondition"
$ws.Range("E1").Value = "Comments"

# Row 2
$ws.Range("A2").Value = "FailToClearTableWrongCode"
$ws.Range("B2").Value = "Costumer ask for the bill or log in through the terminal`nSystem procure bill for the costumer based on his status Subscriber or Guest`nCostumer insert WRONG confirmation code`nSystem procure message `"Wrong confirmation code`""
$ws.Range("B2").WrapText = $true
$ws.Range("B2").VerticalAlignment = -4108
$ws.Range("C2").Value = "Current table has not been payed`nwrong confirmation code given`nTable could not be cleared"
$ws.Range("C2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 390

# Row 3
$ws.Range("A3").Value = "FailToClearTableEmptyCode"
$ws.Range("B3").Value = "Costumer ask for the bill or log in through the terminal`nSystem procure bill for the costumer based on his status Subscriber or Guest`nCostumer insert Empty confirmation code`nSystem procure message `"Wrong confirmation code`""
$ws.Range("B3").WrapText = $true
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("C3").Value = "Current table has not been payed`nwrong confirmation code given`nTable could not be cleared"
$ws.Range("C3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 390

# Row 4
$ws.Range("A4").Value = "FailToClearTablePaymentFailed"
$ws.Range("B4").Value = "Costumer ask for the bill or log in through the terminal`nCostumer insert Correct confirmation code`nCostumer tries to pay his bill`nSystem check wether Payment occured"
$ws.Range("B4").WrapText = $true
$ws.Range("B4").VerticalAlignment = -4108
$ws.Range("C4").Value = "Payment failed`nTable could not be cleared"
$ws.Range("C4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 315

# Row 5
$ws.Range("A5").Value = "TableClearedSuccessfullly"
$ws.Range("B5").Value = "Costumer ask for the bill or log in through the terminal`nSystem procure bill for the costumer based on his status Subscriber or Guest`nCostumer insert Correct confirmation code`nCostumer tries to pay his bill`nSystem check wether Payment occurred"
$ws.Range("B5").WrapText = $true
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("C5").Value = "Payment Succeded`nTable Cleared Successfullly"
$ws.Range("C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 409.5

# Row 6
$ws.Range("A6").Value = "BillDiscountTest"
$ws.Range("B6").Value = "Costumer ask for the bill`nSystem procure bill for the costumer based on his status Subscriber or Guest`nSystem recognize it is a Subscriber`nSystem procures bill with 10% discount`nCostumer insert Correct confirmation code`nCostumer tries to pay his bill`nSystem check wether Payment occurred"
$ws.Range("B6").WrapText = $true
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("C6").Value = "Payment Succeded`nTable Cleared Successfullly"
$ws.Range("C6").WrapText = $true
$ws.Range("D6").Value = "Costumer is a subscriber"
$ws.Rows.Item(6).RowHeight = 409.5

# Row 7
$ws.Range("A7").Value = "SystemBillMessageFailCanceledTable"
$ws.Range("B7").Value = "System recognize two hours`nSystem  try to procure bill for the costumer based on his status Subscriber or Guest`nSystem could not create the bill"
$ws.Range("B7").WrapText = $true
$ws.Range("B7").VerticalAlignment = -4108
$ws.Range("C7").Value = "Bill could not be made.`nTable canceled"
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = "Two Hours Passed"
$ws.Rows.Item(7).RowHeight = 240

# Row 8
$ws.Range("A8").Value = "SystemBillMessageSuccess"
$ws.Range("B8").Value = "System recognize two hours`nSystem try to procure bill for the costumer based on his status Subscriber or Guest`nSystem sends bill to costumer"
$ws.Range("B8").WrapText = $true
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("C8").Value = "Bill Sent to costumer"
$ws.Range("C8").WrapText = $true
$ws.Range("D8").Value = "Two Hours Passed"
$ws.Rows.Item(8).RowHeight = 240

# Column widths
$ws.Columns.Item(2).ColumnWidth = 22.72
$ws.Columns.Item(3).ColumnWidth = 25.39
$ws.Columns.Item(4).ColumnWidth = 19.72
$ws.Columns.Item(5).ColumnWidth = 8.61

# Selection matches the authored file (cell D2 selected on the new active sheet)
[void]$ws.Range("D2").Select()